$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 129
$ws.Range("A3").Value = 178.3999999999996
$ws.Range("A4").Value = 136
$ws.Range("A5").Value = 180.3999999999996
$ws.Range("A6").Value = 126
$ws.Range("A7").Value = 181.3999999999996
$ws.Range("A8").Value = 186.3999999999996
$ws.Range("A9").Value = 192.3999999999996
